# Update the "Forecast Comparison" sheet with the corrected forecast output:
#  - insert a new "Week_Start_Date" column between "Week" and "ASIN"
#  - populate it with each week's start date
#  - normalize the week labels from "W01".."W16" to "W1".."W16"
#  - store is_holiday_week as a boolean instead of a 0/1 number

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column before the old column B (ASIN), shifting ASIN..is_holiday_week
# one column to the right (B->C, C->D, ..., I->J).
$ws.Columns.Item(2).Insert()

# New header for the inserted column.
$ws.Range("B1").Value = "Week_Start_Date"

# Week-start dates (Sundays) for weeks 1-16, aligned with rows 2-17.
$weekStartDates = @(
    "2025-01-05", "2025-01-12", "2025-01-19", "2025-01-26",
    "2025-02-02", "2025-02-09", "2025-02-16", "2025-02-23",
    "2025-03-02", "2025-03-09", "2025-03-16", "2025-03-23",
    "2025-03-30", "2025-04-06", "2025-04-13", "2025-04-20"
)

# Format the new column as text first so the date-like strings are stored
# verbatim (e.g. "2025-01-05") instead of being auto-converted to date serials.
$ws.Range("B2:B17").NumberFormat = "@"

for ($i = 0; $i -lt $weekStartDates.Length; $i++) {
    $row = $i + 2

    # Write the week start date into the newly inserted column B.
    $ws.Range("B$row").Value = $weekStartDates[$i]

    # Normalize "W01" -> "W1", "W02" -> "W2", ... "W16" -> "W16".
    $weekLabel = $ws.Range("A$row").Value2
    $weekNumber = [int]($weekLabel -replace '^W0*', '')
    $ws.Range("A$row").Value = "W$weekNumber"
}

# is_holiday_week now lives in column J (shifted from I) and should be a
# proper boolean rather than a numeric 0/1.
for ($row = 2; $row -le 17; $row++) {
    $flag = $ws.Range("J$row").Value2
    $ws.Range("J$row").Value = [bool]$flag
}
